$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.872.58"
$ws.Range("E2").Value = "  -0.73%  "
$ws.Range("D3").Value = "1.860.92"
$ws.Range("E3").Value = "  -0.33%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.76"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5059"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3633"
$ws.Range("D8").NumberFormat = "General"
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07175"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8953"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  +1.01%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.69"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +0.46%  "
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07474"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = "  -1.04%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.853.08"
$ws.Range("E13").Value = "  -0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.49"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +3.60%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.229"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = "  -1.62%  "
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008473"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +0.00%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.15"
$ws.Range("D18").NumberFormat = "General"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  -0.06%  "
$ws.Range("D20").Value = "26.912.84"
$ws.Range("E20").Value = "  -0.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.026"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  -0.62%  "
$ws.Range("D22").Value = "2.090.51"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.36"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.406"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  -1.26%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "148.01"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  -1.78%  "
$ws.Range("E26").Value = "  -2.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.87"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.055"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -2.15%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "113.06"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.680"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.674"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  -0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09259"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +2.42%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05089"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7423"
$ws.Range("D34").NumberFormat = "General"
$ws.Range("E34").Value = "  +0.60%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.961"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  -4.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.147"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  -0.96%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.280"
$ws.Range("D37").NumberFormat = "General"
$ws.Range("E37").Value = "  +7.80%  "
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.500"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5531"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = "  +3.69%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.070"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  -1.01%  "
$ws.Range("E42").Value = "  +2.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.474"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = "  -1.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.485"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  +1.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1469"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  -0.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4688"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").Value = "  +1.11%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.9997"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  -0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.993"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").Value = "  -0.22%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.564"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.99"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +1.70%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.94"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  -2.55%  "
